$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(160, 1).Value = "2023-12-10 14:07:33"
$ws.Cells.Item(160, 2).Value = 0.0006000000000000001

$ws.Cells.Item(161, 1).Value = "2023-12-10 14:07:49"
$ws.Cells.Item(161, 2).Value = 0.0004
